$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Public Health" column (H) values for rows 2-42: replace the
# placeholder text "---" with the tiny numeric value used by the author.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 8).Value = 0.00000000001
}

# Update the selected cell/range shown in the saved view.
$ws.Range("F5").Select()
